$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Vanguard")

# --- Data updates on the Vanguard scene/faction grid (rows 21-42) ---

# Row 21-22: alpha1 / alpha 2 & sigma1 / sigma 2 already present, unchanged.

# Row 24 (H1): add zeta appearing in column E.
$ws.Range("E24").Value = "zeta"

# Row 25 (H2): column D swaps to "beta2", add zeta in column E.
$ws.Range("D25").Value = "beta2"
$ws.Range("E25").Value = "zeta"

# Row 27 (V5): alpha3 -> alpha2/alpha3, sigma2/sigma3 -> sigma1/sigma2
$ws.Range("B27").Value = "alpha2 / alpha3"
$ws.Range("C27").Value = "sigma1 / sigma2"

# Row 28 (V6): add delta1 in column F.
$ws.Range("F28").Value = "delta1"

# Row 29 (V7): full realignment.
$ws.Range("B29").Value = "alpha3"
$ws.Range("C29").Value = "sigma2 / sigma3"
$ws.Range("D29").Value = "beta2 / beta3 / gamma1 / gamma2"
$ws.Range("E29").Value = "delta1"
$ws.Range("F29").Value = ""

# Row 31 (V8): beta2/beta3 -> beta1/beta2
$ws.Range("D31").Value = "beta1 / beta2"

# Row 32 (V9): swap gamma1/gamma2 <-> sigma2/sigma3 between C and E.
$ws.Range("C32").Value = "sigma2 / sigma3"
$ws.Range("E32").Value = "gamma1 / gamma2"

# Row 37 (H5): drop sigma2/sigma3, shift gamma3/delta2 up.
$ws.Range("D37").Value = "gamma3"
$ws.Range("E37").Value = "delta2"
$ws.Range("F37").Value = ""

# Row 38 (H6): drop sigma2/sigma3, shift gamma.../delta3 up.
$ws.Range("D38").Value = "gamma1 / gamma2 / gamma3"
$ws.Range("E38").Value = "delta3"
$ws.Range("F38").Value = ""

# Row 40 (V10): alpha2/sigma2 -> alpha2/alpha3, sigma2/sigma3 normalized spacing, add zeta in G.
$ws.Range("B40").Value = "alpha2 / alpha3"
$ws.Range("D40").Value = "sigma2 / sigma3"
$ws.Range("G40").Value = "zeta"

# Row 41 (V11): alpha3/sigma3 -> alpha2/alpha3, sigma2/sigma3 normalized spacing, add zeta in G.
$ws.Range("B41").Value = "alpha2 / alpha3"
$ws.Range("D41").Value = "sigma2 / sigma3"
$ws.Range("G41").Value = "zeta"

# Row 42 (V12): alpha3/sigma3 -> alpha2/alpha3, sigma2/sigma3 normalized spacing.
$ws.Range("B42").Value = "alpha2 / alpha3"
$ws.Range("D42").Value = "sigma2 / sigma3"

# --- Column D widened to fit the longer "beta2 / beta3 / gamma1 / gamma2" text ---
$ws.Columns.Item(4).ColumnWidth = 30.8

# --- Update the active selection on the Vanguard sheet ---
$ws.Range("H37").Select()
